$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3 (PORCELANATO)
$ws.Range("D3").Value = 1331.8
$ws.Range("E3").Value = 12391.54
$ws.Range("F3").Value = 0.09704634586041007

# Row 4 (TOTAL)
$ws.Range("D4").Value = 4924.47
$ws.Range("E4").Value = 8798.870000000001
$ws.Range("F4").Value = 0.3588390289827404
